$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellAddress,
        [string]$Text
    )
    # Assigning certain strings (e.g. "34%") directly to a Range.Value causes Excel
    # to auto-recognize them as numbers/percentages and reformat the cell.
    # Routing the literal text through a formula and then collapsing it to a
    # static value via Copy / PasteSpecial (xlPasteValues) keeps the cell's
    # content as plain text without attaching any extra number format/style.
    $escaped = $Text.Replace('"', '""')
    $range = $ws.Range($CellAddress)
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# --- Row 28 ---
Set-TextValue "A28" "2024-10-10 21:29:36"
$ws.Range("C28").Value = 33
$ws.Range("D28").Value = 10
$ws.Range("E28").Value = 23
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("N28").Value = 10
$ws.Range("O28").Value = 10
$ws.Range("P28").Value = 1
$ws.Range("R28").Value = 5
$ws.Range("T28").Value = 20
Set-TextValue "U28" "34%"
Set-TextValue "V28" "C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Crupier.xlsx"
Set-TextValue "X28" "No es Simulación"
$ws.Range("Y28").Value = 98

# --- Row 29 ---
Set-TextValue "A29" "2024-10-10 21:58:24"
$ws.Range("C29").Value = 16
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = 9
$ws.Range("F29").Value = 5
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("N29").Value = 10
$ws.Range("O29").Value = 10
$ws.Range("P29").Value = 2
$ws.Range("R29").Value = 5
$ws.Range("T29").Value = 20
Set-TextValue "U29" "47%"
Set-TextValue "V29" "C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Crupier.xlsx"
Set-TextValue "X29" "No es Simulación"
$ws.Range("Y29").Value = 34
